$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (16-core results block, columns F:H) ---
$ws.Range("G1").Value = "16 ядер"

$ws.Range("F2").Value = "До"
$ws.Range("G2").Value = "После"
$ws.Range("H2").Value = "Ускорение"

# --- Row 9: convolutionConcurrentNN / Общее (2 итерации) ---
# (written before B3 so shared-string order matches: 16 ядер, convolutionConcurrentNN, 42 значения)
$ws.Range("A9").Value = "convolutionConcurrentNN"
$ws.Range("B9").Value = "Общее (2 итерации)"
$ws.Range("C9").Value = 1.723
$ws.Range("F9").Value = 17.1018

# --- Row 3: SimpleExample ---
$ws.Range("B3").Value = "42 значения"
$ws.Range("F3").Value = 12.1088
$ws.Range("G3").Value = 1.96376
$ws.Range("H3").Formula = "=F3/G3"

# --- Row 4: ImageScaleParallel ---
$ws.Range("F4").Value = 14.848697
$ws.Range("G4").Value = 4.96525
$ws.Range("H4").Formula = "=F4/G4"

# --- Row 5: convolutionNN / Общее (2 итерации) ---
$ws.Range("F5").Value = 0.8217

# --- Row 6: обучение на итерации ---
$ws.Range("F6").Value = 0.3158

# --- Row 7: тест на итерации ---
$ws.Range("F7").Value = 0.09437

# --- Row 10: обучение на итерации ---
$ws.Range("B10").Value = "обучение на итерации"
$ws.Range("C10").Value = 0.512
$ws.Range("F10").Value = 0.3118

# --- Row 11: тест на итерации ---
$ws.Range("B11").Value = "тест на итерации"
$ws.Range("C11").Value = 0.343
$ws.Range("F11").Value = 8.2147

# --- Column A best-fit width (new longest label is "convolutionConcurrentNN") ---
$ws.Columns.Item(1).EntireColumn.AutoFit()

# --- Selection moves to B2 ---
[void]$ws.Range("B2").Select()
